$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Value = "test"
$ws.Range("L7").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
Write-Host "done"
